# Removed white space from Fig 1
# - Group the loose oval/number callouts, the base picture, the nested
#   "Group 1" callout and the small cropped picture into a single new
#   group shape (matches PowerPoint's "Group" command which was used to
#   tidy up Figure 1 on this slide).
# - Fix a typo in the table: "3b" -> "3a".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- 1. Fix the typo in the table cell text -----------------------------
$tbl = $s.Shapes.Item("Table 7").Table
$cell = $tbl.Cell(9, 2)
$cell.Shape.TextFrame.TextRange.Text = "Reduces value of 3a by 50%, or 90%"

# --- 2. Group every shape on the slide except the table into one group --
$namesToGroup = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -ne "Table 7") {
        $namesToGroup += $sh.Name
    }
}

$range = $s.Shapes.Range($namesToGroup)
$newGroup = $range.Group()
$newGroup.Name = "Group 2"
